# "navn for dag 1" - fill in the runner names for Day 1 heats (K1, K2, K3)
$wb = $excel.ActiveWorkbook

$k1 = $wb.Worksheets.Item("K1")
$k2 = $wb.Worksheets.Item("K2")
$k3 = $wb.Worksheets.Item("K3")

# Order matters here: it reproduces the exact shared-string insertion
# order from the authored workbook.
$k1.Range("A5").Value = "Jonatan Heimdal Korshavn"
$k1.Range("A7").Value = "Marius Granvold"
$k1.Range("A4").Value = "Magnus Moslett Evensen"
$k1.Range("A6").Value = "Martin Jørstad Ringli"

$k2.Range("A4").Value = "Eskil Engdal"
$k2.Range("A5").Value = "Tobias Gigstad Bergene"

$k3.Range("A4").Value = "Brage Sømoen"
$k3.Range("A5").Value = "Jakob Lundby"
$k3.Range("A6").Value = "Anders Alme Eng"
$k3.Range("A7").Value = "Vegard Thon"
$k3.Range("A8").Value = "Christian Thon Christensen"
$k3.Range("A9").Value = "Per Ingvar Tollehaug"
$k3.Range("A10").Value = "Hermann Byfuglien Ulsrud"
$k3.Range("A11").Value = "Erland Andersen"
